$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.6308291129860543
$ws.Range("C4").Value = 0.657
$ws.Range("D4").Value = 0.624000858260329
$ws.Range("E4").Value = 0.6234999999999999
$ws.Range("F4").Value = 0.5652977737532658
$ws.Range("G4").Value = 0.552
$ws.Range("H4").Value = 0.5977215165268509
$ws.Range("I4").Value = 0.58
$ws.Range("J4").Value = 0.6519849288583972
$ws.Range("K4").Value = 0.843
$ws.Range("L4").Value = 0.5380798606831865
$ws.Range("M4").Value = 0.5565
$ws.Range("N4").Value = 0.6395628362439414
$ws.Range("O4").Value = 0.634
$ws.Range("P4").Value = 0.6627808854146335
$ws.Range("Q4").Value = 0.643
$ws.Range("R4").Value = 0.63549681758873
$ws.Range("S4").Value = 0.654
$ws.Range("T4").Value = 0.6395795601578353
$ws.Range("U4").Value = 0.6345000000000001
$ws.Range("V4").Value = 0.5649906893790165
$ws.Range("W4").Value = 0.5569999999999999
$ws.Range("X4").Value = 0.5882215793496119
$ws.Range("Y4").Value = 0.5745
$ws.Range("Z4").Value = 0.6442170220377668
$ws.Range("AA4").Value = 0.6439999999999999
$ws.Range("AB4").Value = 0.6611430738450725
$ws.Range("AC4").Value = 0.6439999999999999

$ws.Range("B5").Value = 0.6424736233774221
$ws.Range("C5").Value = 0.694
$ws.Range("D5").Value = 0.6105672132414932
$ws.Range("E5").Value = 0.625
$ws.Range("F5").Value = 0.5984917585207542
$ws.Range("G5").Value = 0.593
$ws.Range("H5").Value = 0.6228218685650765
$ws.Range("I5").Value = 0.6035000000000001
$ws.Range("J5").Value = 0.6037038077613801
$ws.Range("K5").Value = 0.724
$ws.Range("L5").Value = 0.5379406927739231
$ws.Range("M5").Value = 0.553
$ws.Range("N5").Value = 0.6161665065075501
$ws.Range("O5").Value = 0.632
$ws.Range("P5").Value = 0.6217610241101575
$ws.Range("Q5").Value = 0.609
$ws.Range("R5").Value = 0.6337188302597639
$ws.Range("S5").Value = 0.6839999999999999
$ws.Range("T5").Value = 0.6064304464508008
$ws.Range("U5").Value = 0.6199999999999999
$ws.Range("V5").Value = 0.5990662892525185
$ws.Range("W5").Value = 0.594
$ws.Range("X5").Value = 0.6229116812850777
$ws.Range("Y5").Value = 0.6035
$ws.Range("Z5").Value = 0.6155354443603589
$ws.Range("AA5").Value = 0.631
$ws.Range("AB5").Value = 0.6214260710346802
$ws.Range("AC5").Value = 0.608

$ws.Range("B6").Value = 0.644080694845645
$ws.Range("C6").Value = 0.658
$ws.Range("D6").Value = 0.6486519887610067
$ws.Range("E6").Value = 0.6465
$ws.Range("F6").Value = 0.5996900752419304
$ws.Range("G6").Value = 0.6220000000000001
$ws.Range("H6").Value = 0.5936780405527909
$ws.Range("I6").Value = 0.5894999999999999
$ws.Range("J6").Value = 0.6457397050733061
$ws.Range("K6").Value = 0.833
$ws.Range("L6").Value = 0.53382230517117
$ws.Range("M6").Value = 0.5505
$ws.Range("N6").Value = 0.6242481815970227
$ws.Range("O6").Value = 0.634
$ws.Range("P6").Value = 0.6338621883515176
$ws.Range("Q6").Value = 0.6204999999999999
$ws.Range("R6").Value = 0.6456099236294124
$ws.Range("S6").Value = 0.659
$ws.Range("T6").Value = 0.6520829127344496
$ws.Range("U6").Value = 0.648
$ws.Range("V6").Value = 0.6036240071995228
$ws.Range("W6").Value = 0.628
$ws.Range("X6").Value = 0.5952009161514079
$ws.Range("Y6").Value = 0.592
$ws.Range("Z6").Value = 0.6269390047868134
$ws.Range("AA6").Value = 0.6380000000000001
$ws.Range("AB6").Value = 0.6347722377263799
$ws.Range("AC6").Value = 0.6219999999999999

$ws.Range("B7").Value = 0.5599715909145864
$ws.Range("C7").Value = 0.574
$ws.Range("D7").Value = 0.5578194224204991
$ws.Range("E7").Value = 0.5525
$ws.Range("F7").Value = 0.5090029130676117
$ws.Range("G7").Value = 0.5229999999999999
$ws.Range("H7").Value = 0.5100693189288688
$ws.Range("I7").Value = 0.505
$ws.Range("J7").Value = 0.5030415363176824
$ws.Range("K7").Value = 0.516
$ws.Range("L7").Value = 0.5091497254670135
$ws.Range("M7").Value = 0.514
$ws.Range("N7").Value = 0.5370677210419921
$ws.Range("O7").Value = 0.5570000000000001
$ws.Range("P7").Value = 0.549942074640928
$ws.Range("Q7").Value = 0.5399999999999999
$ws.Range("R7").Value = 0.5486553910940468
$ws.Range("S7").Value = 0.577
$ws.Range("T7").Value = 0.5369007085352122
$ws.Range("U7").Value = 0.539
$ws.Range("V7").Value = 0.5184461646050238
$ws.Range("W7").Value = 0.5349999999999999
$ws.Range("X7").Value = 0.51703183674548
$ws.Range("Y7").Value = 0.5110000000000001
$ws.Range("Z7").Value = 0.5428169213100973
$ws.Range("AA7").Value = 0.5619999999999999
$ws.Range("AB7").Value = 0.5544563249047233
$ws.Range("AC7").Value = 0.545
